$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the preparation text (F2): Kode Transaksi changes from 090 to 998
$ws.Range("F2").Value = "Username : 30711;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nKode Transaksi : 998"

# Update KODE_JENIS_TRANSAKSI (M2) from text "090" to numeric 998,
# while preserving the cell's existing number format / quote-prefix style.
$ws.Range("M2").Value = 998
$ws.Range("K2").Copy()
$ws.Range("M2").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection left in the sheet when saved (was S2, now G2).
# (The workbook was also scrolled so column D is the leftmost visible column;
# the headless engine does not persist plain scroll position to topLeftCell,
# only the active selection, which is applied below.)
$ws.Range("G2").Select()
